# "tambah hazwan dan kholid" — add two new students (Hazwan, Kholid) to the
# nilai ("nilai") sheet, then re-sort the A3:C13 data range alphabetically by
# name (column A), matching the workbook's existing sortState/autosort setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nilai")

# The table previously ended at row 11 (Zufar); rows 12/13 were blank.
# Append the two new students' raw data there first, then sort A3:C13 so the
# whole roster ends up alphabetised, same as the rest of the sheet.
$ws.Range("A12").Value = "Hazwan"
$ws.Range("A13").Value = "Kholid"
$ws.Range("B12").Value = "9 Agustus 2023"
$ws.Range("B13").Value = "9 Agustus 2023"
$ws.Range("C12").Value = 38.9
$ws.Range("C13").Value = 89.2

# Re-sort the roster (A3:C13) ascending by Nama (column A), no header row,
# using the worksheet's Sort object so the persisted sortState/sortCondition
# range is updated to match (A3:C13 / A3:A13).
$sort = $ws.Sort
$sort.SortFields.Clear()
$null = $sort.SortFields.Add($ws.Range("A3:A13"))
$sort.SetRange($ws.Range("A3:C13"))
$sort.Header = 2
$sort.Apply()

# Leave the active cell where the user last clicked while entering the data.
$null = $ws.Range("C5").Select()
